$d = $word.ActiveDocument

# --- Text content reverts (Banking/Finance template swap) ---------------

# Title
$d.Content.Find.Execute("ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING", $true, $false, $false, $false, $false,
                         $true, 1, $false, "FINANCE - CORE BANKING SYSTEM MODERNIZATION", 2)

# "Digital transformation through ..." phrase - occurs 3x (subtitle, exec
# summary sentence, business-context line) and is always a verbatim match.
$d.Content.Find.Execute("Digital transformation through intelligent automation and predictive analytics", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Digital transformation through modern banking operations and transaction processing", 2)

# "Finance and Machine Learning" - occurs 4x (Industry, exec summary
# sentence, Project Name, Industry Focus).
$d.Content.Find.Execute("Finance and Machine Learning", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Banking and Banking Operations", 2)

# "Finance Implementation" - occurs 5x (Project Type x2, exec summary
# sentence, strategic alignment sentence, proposed solution sentence).
$d.Content.Find.Execute("Finance Implementation", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Banking Implementation", 2)

# --- Add a (blank) default header/footer pair to the only section -------

$sec = $d.Sections(1)

$hdr = $sec.Headers(1)          # wdHeaderFooterPrimary
$hdr.Range.Paragraphs(1).Style = "Header"

$ftr = $sec.Footers(1)          # wdHeaderFooterPrimary
$ftr.Range.Paragraphs(1).Style = "Footer"
